$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row: C1="prediction", D1="rejection-f", E1="max"
$ws.Range("C1").Value = "prediction"
$ws.Range("D1").Value = "rejection-f"
$ws.Range("E1").Value = "max"

# New numeric values for column E (was text duplicated from D, now a numeric score)
$eValues = @(
    0.9646825789224177,
    0.968021397348124,
    0.9698807119560511,
    0.9691163253920958,
    0.9684825329576824,
    0.9698742099727952,
    0.9692630485448526,
    0.969752441022579,
    0.9694952613715876,
    0.9688388989897893,
    0.9685586922578093,
    0.9693606522736604,
    0.9704040846133849,
    0.9686882616868068,
    0.965900285551642,
    0.9700699439746875,
    0.9703082058493662,
    0.968480033605666,
    0.9700754735150953
)

# Column C previously duplicated column B's numeric value; it now duplicates
# column D's species text instead. Column E previously duplicated column D's
# text too; it now holds its own numeric value.
for ($i = 0; $i -lt $eValues.Length; $i++) {
    $row = $i + 2
    $species = $ws.Range("D$row").Value()
    $ws.Range("C$row").Value = $species
    $ws.Range("E$row").Value = $eValues[$i]
}
